# Add a new "created_at" date column (F) to the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("F1").Value = "created_at"

# Data (serial date values: 2024-02-03 .. 2024-02-06)
$ws.Range("F2").Value = 45325
$ws.Range("F3").Value = 45326
$ws.Range("F4").Value = 45327
$ws.Range("F5").Value = 45328

# Apply a short-date display format to the new column, then propagate the
# same style to the rest of the rows via a format-only paste so every cell
# shares one cell-format record (instead of one-per-cell).
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F3:F5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the recorded selection after the edit.
$ws.Range("G2").Select() | Out-Null
